$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 236
$ws.Range("I2").Value = 577
$ws.Range("J2").Value = 2449
$ws.Range("K2").Value = 20
$ws.Range("L2").Value = 641
$ws.Range("M2").Value = 33
$ws.Range("N2").Value = 451
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = 7
$ws.Range("R2").Value = 30
$ws.Range("S2").Value = 244
$ws.Range("T2").Value = 415
$ws.Range("U2").Value = 31
$ws.Range("V2").Value = 3799
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 3692
$ws.Range("Y2").Value = 6
$ws.Range("Z2").Value = 55
$ws.Range("AA2").Value = 24
